# Applies the scheduled-runner value refresh to the Leve profit columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1901
$ws.Range("I32").Value = 1800
$ws.Range("J32").Value = 2002
$ws.Range("K32").Value = 1800
$ws.Range("L32").Value = 2002
$ws.Range("M32").Value = -1474
$ws.Range("N32").Value = -2654
$ws.Range("H132").Value = 5994.28
$ws.Range("I132").Value = 8238.764999999999
$ws.Range("J132").Value = 1224.75
$ws.Range("K132").Value = 24716.295
$ws.Range("L132").Value = 3674.25
$ws.Range("M132").Value = -22186.295
$ws.Range("N132").Value = -8734.25
$ws.Range("H138").Value = 3422.5217
$ws.Range("I138").Value = 5377
$ws.Range("J138").Value = 3184.1707
$ws.Range("K138").Value = 16131
$ws.Range("L138").Value = 9552.5121
$ws.Range("M138").Value = -10991
$ws.Range("N138").Value = -19832.5121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31935.475
$ws.Range("I32").Value = 6113.4634
$ws.Range("J32").Value = 98104.375
$ws.Range("K32").Value = 6113.4634
$ws.Range("L32").Value = 98104.375
$ws.Range("M32").Value = -5826.4634
$ws.Range("N32").Value = -98678.375
$ws.Range("H45").Value = 1212.4166
$ws.Range("I45").Value = 841.1111
$ws.Range("J45").Value = 2326.3333
$ws.Range("K45").Value = 841.1111
$ws.Range("L45").Value = 2326.3333
$ws.Range("M45").Value = -464.1111
$ws.Range("N45").Value = -3080.3333
$ws.Range("H74").Value = 2203.8667
$ws.Range("I74").Value = 1737.2727
$ws.Range("J74").Value = 2474
$ws.Range("K74").Value = 1737.2727
$ws.Range("L74").Value = 2474
$ws.Range("M74").Value = -863.2727
$ws.Range("N74").Value = -4222
$ws.Range("H77").Value = 2203.8667
$ws.Range("I77").Value = 1737.2727
$ws.Range("J77").Value = 2474
$ws.Range("K77").Value = 8686.363499999999
$ws.Range("L77").Value = 12370
$ws.Range("M77").Value = -4318.363499999999
$ws.Range("N77").Value = -21106
$ws.Range("H132").Value = 4391.7407
$ws.Range("I132").Value = 4716.636
$ws.Range("J132").Value = 2962.2
$ws.Range("K132").Value = 14149.908
$ws.Range("L132").Value = 8886.599999999999
$ws.Range("M132").Value = -11619.908
$ws.Range("N132").Value = -13946.6
$ws.Range("H138").Value = 22517.5
$ws.Range("I138").Value = 22517.5
$ws.Range("K138").Value = 22517.5
$ws.Range("M138").Value = -17377.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 347.4
$ws.Range("I64").Value = 100
$ws.Range("J64").Value = 512.3333
$ws.Range("K64").Value = 100
$ws.Range("L64").Value = 512.3333
$ws.Range("M64").Value = 125
$ws.Range("N64").Value = -962.3333
$ws.Range("H67").Value = 347.4
$ws.Range("I67").Value = 100
$ws.Range("J67").Value = 512.3333
$ws.Range("K67").Value = 100
$ws.Range("L67").Value = 512.3333
$ws.Range("M67").Value = 680
$ws.Range("N67").Value = -2072.3333
$ws.Range("H107").Value = 30342412
$ws.Range("I107").Value = 47677776
$ws.Range("K107").Value = 47677776
$ws.Range("M107").Value = -47675856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 492.16666
$ws.Range("I22").Value = 429.14285
$ws.Range("J22").Value = 580.4
$ws.Range("K22").Value = 429.14285
$ws.Range("L22").Value = 580.4
$ws.Range("M22").Value = -79.14285000000001
$ws.Range("N22").Value = -1280.4
$ws.Range("H31").Value = 1363.9296
$ws.Range("I31").Value = 914.75
$ws.Range("J31").Value = 2301.348
$ws.Range("K31").Value = 914.75
$ws.Range("L31").Value = 2301.348
$ws.Range("M31").Value = -619.75
$ws.Range("N31").Value = -2891.348
$ws.Range("H34").Value = 1363.9296
$ws.Range("I34").Value = 914.75
$ws.Range("J34").Value = 2301.348
$ws.Range("K34").Value = 914.75
$ws.Range("L34").Value = 2301.348
$ws.Range("M34").Value = -712.75
$ws.Range("N34").Value = -2705.348
$ws.Range("H41").Value = 11427.692
$ws.Range("I41").Value = 500
$ws.Range("J41").Value = 12338.333
$ws.Range("K41").Value = 500
$ws.Range("L41").Value = 12338.333
$ws.Range("M41").Value = -72
$ws.Range("N41").Value = -13194.333
$ws.Range("H50").Value = 16940
$ws.Range("J50").Value = 16940
$ws.Range("L50").Value = 16940
$ws.Range("N50").Value = -18190
$ws.Range("H51").Value = 7424.875
$ws.Range("I51").Value = 500
$ws.Range("J51").Value = 7886.533
$ws.Range("K51").Value = 500
$ws.Range("L51").Value = 7886.533
$ws.Range("M51").Value = 236
$ws.Range("N51").Value = -9358.532999999999
$ws.Range("H59").Value = 24845
$ws.Range("J59").Value = 24845
$ws.Range("L59").Value = 24845
$ws.Range("N59").Value = -27135
$ws.Range("H60").Value = 17920
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 17920
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 17920
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -18942
$ws.Range("H61").Value = 7424.875
$ws.Range("I61").Value = 500
$ws.Range("J61").Value = 7886.533
$ws.Range("K61").Value = 500
$ws.Range("L61").Value = 7886.533
$ws.Range("M61").Value = -152
$ws.Range("N61").Value = -8582.532999999999
$ws.Range("H69").Value = 9998
$ws.Range("I69").Value = 9998
$ws.Range("K69").Value = 9998
$ws.Range("M69").Value = -9249
$ws.Range("H72").Value = 9998
$ws.Range("I72").Value = 9998
$ws.Range("K72").Value = 29994
$ws.Range("M72").Value = -26250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1613.5
$ws.Range("I34").Value = 302
$ws.Range("K34").Value = 906
$ws.Range("M34").Value = -822
$ws.Range("H129").Value = 10870916
$ws.Range("I129").Value = 31250456
$ws.Range("J129").Value = 1828.7333
$ws.Range("K129").Value = 93751368
$ws.Range("L129").Value = 5486.199900000001
$ws.Range("M129").Value = -93746368
$ws.Range("N129").Value = -15486.1999
$ws.Range("H131").Value = 901811.0600000001
$ws.Range("J131").Value = 901811.0600000001
$ws.Range("L131").Value = 2705433.18
$ws.Range("N131").Value = -2715513.18
$ws.Range("H133").Value = 2932.76
$ws.Range("I133").Value = 634.6
$ws.Range("K133").Value = 1903.8
$ws.Range("M133").Value = 3156.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 204190.3
$ws.Range("I70").Value = 289700.72
$ws.Range("J70").Value = 4666
$ws.Range("K70").Value = 289700.72
$ws.Range("L70").Value = 4666
$ws.Range("M70").Value = -289430.72
$ws.Range("N70").Value = -5206
$ws.Range("H73").Value = 204190.3
$ws.Range("I73").Value = 289700.72
$ws.Range("J73").Value = 4666
$ws.Range("K73").Value = 289700.72
$ws.Range("L73").Value = 4666
$ws.Range("M73").Value = -288764.72
$ws.Range("N73").Value = -6538
$ws.Range("H132").Value = 2169.4
$ws.Range("I132").Value = 1899.1578
$ws.Range("J132").Value = 2636.182
$ws.Range("K132").Value = 5697.4734
$ws.Range("L132").Value = 7908.545999999999
$ws.Range("M132").Value = -3167.4734
$ws.Range("N132").Value = -12968.546
$ws.Range("H135").Value = 40417.332
$ws.Range("J135").Value = 40417.332
$ws.Range("L135").Value = 40417.332
$ws.Range("N135").Value = -50557.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 79977.234
$ws.Range("I40").Value = 252501
$ws.Range("J40").Value = 3300
$ws.Range("K40").Value = 252501
$ws.Range("L40").Value = 3300
$ws.Range("M40").Value = -252365
$ws.Range("N40").Value = -3572
$ws.Range("H46").Value = 1265974
$ws.Range("J46").Value = 1265974
$ws.Range("L46").Value = 1265974
$ws.Range("N46").Value = -1266350
$ws.Range("H122").Value = 3325
$ws.Range("I122").Value = 3325
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9975
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7525
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1235.9546
$ws.Range("I136").Value = 642.25
$ws.Range("J136").Value = 1948.4
$ws.Range("K136").Value = 1926.75
$ws.Range("L136").Value = 5845.200000000001
$ws.Range("M136").Value = 623.25
$ws.Range("N136").Value = -10945.2
$ws.Range("H140").Value = 35195
$ws.Range("I140").Value = 6390
$ws.Range("J140").Value = 64000
$ws.Range("K140").Value = 6390
$ws.Range("L140").Value = 64000
$ws.Range("M140").Value = -1210
$ws.Range("N140").Value = -74360
